$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4809054136276245
$ws.Range("B1").Value = 0.3935800194740295
$ws.Range("C1").Value = 0.3510326743125916
$ws.Range("D1").Value = 0.3866152763366699
$ws.Range("E1").Value = 0.482348769903183
